$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 286, pushing the existing rows 286-400 down to 287-401
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row 286 with the new record
$ws.Cells.Item(286, 1).Value  = 10
$ws.Cells.Item(286, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(286, 3).Value  = 'La Araucanía'
$ws.Cells.Item(286, 4).Value  = 45009
$ws.Cells.Item(286, 5).Value  = 9
$ws.Cells.Item(286, 6).Value  = 'Fruta'
$ws.Cells.Item(286, 7).Value  = 100102
$ws.Cells.Item(286, 8).Value  = 'Cítricos'
$ws.Cells.Item(286, 9).Value  = 100102006
$ws.Cells.Item(286, 10).Value = 'Pomelo'
$ws.Cells.Item(286, 11).Value = 'Start Ruby'
$ws.Cells.Item(286, 12).Value = 'Primera'
$ws.Cells.Item(286, 13).Value = 80
$ws.Cells.Item(286, 14).Value = 14000
$ws.Cells.Item(286, 15).Value = 14000
$ws.Cells.Item(286, 16).Value = 14000
$ws.Cells.Item(286, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(286, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(286, 19).Value = 933
$ws.Cells.Item(286, 20).Value = 15
